$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.617.33'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.638.53'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.10%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '512.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.89'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.91%  '
$ws.Range("E7").Value = '  -0.35%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.565'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.65%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.667.31'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.24'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.57%  '
$ws.Range("E11").Value = '  +3.58%  '
$ws.Range("E12").Value = '  +0.51%  '
$ws.Range("E13").Value = '  -1.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.103.91'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.23%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '58.619.54'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.07%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.84'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.11%  '
$ws.Range("E17").Value = '  +1.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.663.56'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.52'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '342.82'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.74%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.36'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.78%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.09'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.37%  '
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.64'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.46%  '
$ws.Range("E25").Value = '  +2.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.761.76'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.993'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.68%  '
$ws.Range("E28").Value = '  +1.72%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0802'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.16'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.34'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +8.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.80'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.57'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '149.48'
$ws.Range("D35").Style = "Normal"
$ws.Range("E36").Value = '  +10.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.98'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.35%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.14'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.75%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.57'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.842'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.46%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.66'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.72%  '
$ws.Range("E42").Value = '  +0.75%  '
$ws.Range("E43").Value = '  +1.54%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '277.35'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.994'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.37%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0976'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.47'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.24%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0529'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.77%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0229'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.15%  '
$ws.Range("B50").Value = 'WhiteBITCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '10.26'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.55%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.69'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.32%  '
